$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are numeric-looking text (e.g. thousand-dot-grouped
# prices like "34.011.08"); force text format so Excel keeps them as strings
# instead of parsing them into floating point numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.011.08"
$ws.Range("E2").Value = "  -1.70%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.791.10"
$ws.Range("E3").Value = "  -1.88%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.96"
$ws.Range("E5").Value = "  -2.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.553"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.10"
$ws.Range("E8").Value = "  -3.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.01"
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.280"
$ws.Range("E10").Value = "  -2.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0660"
$ws.Range("E11").Value = "  -3.37%  "
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.050.86"
$ws.Range("E13").Value = "  -1.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.34"
$ws.Range("E14").Value = "  +8.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.793.08"
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.633"
$ws.Range("E16").Value = "  -2.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "34.053.07"
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.21"
$ws.Range("E18").Value = "  -2.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.43"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "252.70"
$ws.Range("E20").Value = "  -3.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0744"
$ws.Range("E21").Value = "  -2.44%  "
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.41"
$ws.Range("E23").Value = "  -1.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.28"
$ws.Range("E24").Value = "  -3.38%  "
$ws.Range("E25").Value = "  -2.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.52"
$ws.Range("E26").Value = "  -2.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.59"
$ws.Range("E27").Value = "  -2.46%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.01"
$ws.Range("E28").Value = "  -2.71%  "
$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.114"
$ws.Range("E29").Value = "  -2.87%  "
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.88"
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.62"
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.84"
$ws.Range("E35").Value = "  -1.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.474.48"
$ws.Range("E36").Value = "  -7.25%  "
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.632"
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0186"
$ws.Range("E39").Value = "  -1.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "83.89"
$ws.Range("E40").Value = "  -3.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.82"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.901"
$ws.Range("E43").Value = "  -2.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.05"
$ws.Range("E44").Value = "  -4.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0515"
$ws.Range("E45").Value = "  -2.61%  "
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.950.14"
$ws.Range("E47").Value = "  -1.44%  "
$ws.Range("E48").Value = "  +0.31%  "
$ws.Range("E49").Value = "  -1.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.76"
$ws.Range("E50").Value = "  +2.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.26"
$ws.Range("E51").Value = "  -5.69%  "
